# CUS15: ajustes en controlador y cotizaciones
# Update the quotation line items (products, quantities and ids) on the
# "COTIZACION" sheet. The totals (E17:E19, E21:E23) are driven by existing
# formulas and will recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COTIZACION")

# Row 19: Id_Producto / Descripcion / CantidadOfertada
$ws.Range("A19").Value = 1028
$ws.Range("B19").Value = "Dentastix para perros pequeños 3pzas"
$ws.Range("C19").Value = 200

# Row 18: Id_Producto / Descripcion / CantidadOfertada
$ws.Range("A18").Value = 1020
$ws.Range("B18").Value = "Hueso Dental Chiquito"
$ws.Range("C18").Value = 150

# Row 17: Id_Producto / Descripcion / CantidadOfertada
$ws.Range("A17").Value = 1018
$ws.Range("B17").Value = "Felix Party Mix 60g"
$ws.Range("C17").Value = 200

$excel.Calculate()

# Restore scroll position / selection as recorded in the view state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B21").Select()
